# Add team record (Wins/Losses/Ties) columns to the data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the bold/bordered header style used by the existing header cells (e.g. AC1)
# by copying the formatting of an existing header cell onto the new ones.
$headerFormat = $ws.Range("AC1")
$headerFormat.Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Data rows 2-45: every row gets the same team record (88-74-0)
for ($row = 2; $row -le 45; $row++) {
    $ws.Cells.Item($row, 30).Value = 88   # AD -> Wins
    $ws.Cells.Item($row, 31).Value = 74   # AE -> Losses
    $ws.Cells.Item($row, 32).Value = 0    # AF -> Ties
}
